# Auto-generated: applies the 'Update countries & provincias Spain' edit
# to paises.xlsx - refreshes COVID-19 country stats and the as-of timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 4 de Octubre de 2020 a las 16:21'
$ws.Cells.Item(4, 2).Value = 7603074
$ws.Cells.Item(4, 3).Value = 2228
$ws.Cells.Item(4, 4).Value = 4818788
$ws.Cells.Item(4, 5).Value = 2570000
$ws.Cells.Item(4, 7).Value = 9
$ws.Cells.Item(4, 8).Value = 214286
$ws.Cells.Item(5, 2).Value = 6559777
$ws.Cells.Item(5, 3).Value = 12364
$ws.Cells.Item(5, 4).Value = 5516487
$ws.Cells.Item(5, 5).Value = 941380
$ws.Cells.Item(5, 7).Value = 98
$ws.Cells.Item(5, 8).Value = 101910
$ws.Cells.Item(17, 2).Value = 470179
$ws.Cells.Item(17, 3).Value = 1708
$ws.Cells.Item(17, 4).Value = 442070
$ws.Cells.Item(17, 5).Value = 15130
$ws.Cells.Item(17, 7).Value = 60
$ws.Cells.Item(17, 8).Value = 12979
$ws.Cells.Item(40, 2).Value = 107025
$ws.Cells.Item(40, 3).Value = 567
$ws.Cells.Item(40, 4).Value = 99040
$ws.Cells.Item(40, 5).Value = 7361
$ws.Cells.Item(40, 7).Value = 4
$ws.Cells.Item(40, 8).Value = 624
$ws.Cells.Item(53, 1).Value = 'Portugal'
$ws.Cells.Item(53, 2).Value = 79151
$ws.Cells.Item(53, 3).Value = 904
$ws.Cells.Item(53, 4).Value = 50207
$ws.Cells.Item(53, 5).Value = 26939
$ws.Cells.Item(53, 7).Value = 10
$ws.Cells.Item(53, 8).Value = 2005
$ws.Cells.Item(54, 1).Value = 'Honduras'
$ws.Cells.Item(54, 2).Value = 78788
$ws.Cells.Item(54, 3).Value = 519
$ws.Cells.Item(54, 4).Value = 29187
$ws.Cells.Item(54, 5).Value = 47202
$ws.Cells.Item(54, 7).Value = 13
$ws.Cells.Item(54, 8).Value = 2399
$ws.Cells.Item(60, 4).Value = 57575
$ws.Cells.Item(60, 5).Value = 210
$ws.Cells.Item(72, 1).Value = 'Kenia'
$ws.Cells.Item(72, 2).Value = 39427
$ws.Cells.Item(72, 3).Value = 243
$ws.Cells.Item(72, 4).Value = 25659
$ws.Cells.Item(72, 5).Value = 13037
$ws.Cells.Item(72, 7).Value = 3
$ws.Cells.Item(72, 8).Value = 731
$ws.Cells.Item(73, 1).Value = 'Afganistan'
$ws.Cells.Item(73, 2).Value = 39341
$ws.Cells.Item(73, 3).Value = 44
$ws.Cells.Item(73, 4).Value = 32852
$ws.Cells.Item(73, 5).Value = 5027
$ws.Cells.Item(73, 8).Value = 1462
$ws.Cells.Item(75, 2).Value = 36809
$ws.Cells.Item(75, 3).Value = 722
$ws.Cells.Item(75, 4).Value = 21429
$ws.Cells.Item(75, 5).Value = 14788
$ws.Cells.Item(75, 7).Value = 14
$ws.Cells.Item(75, 8).Value = 592
$ws.Cells.Item(76, 2).Value = 33901
$ws.Cells.Item(76, 3).Value = 59
$ws.Cells.Item(76, 5).Value = 1611
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 754
$ws.Cells.Item(79, 2).Value = 29450
$ws.Cells.Item(79, 3).Value = 92
$ws.Cells.Item(79, 4).Value = 24296
$ws.Cells.Item(79, 5).Value = 4291
$ws.Cells.Item(80, 2).Value = 28354
$ws.Cells.Item(80, 3).Value = 379
$ws.Cells.Item(80, 4).Value = 21859
$ws.Cells.Item(80, 5).Value = 5616
$ws.Cells.Item(80, 7).Value = 8
$ws.Cells.Item(80, 8).Value = 879
$ws.Cells.Item(88, 2).Value = 18790
$ws.Cells.Item(88, 3).Value = 188
$ws.Cells.Item(88, 4).Value = 15341
$ws.Cells.Item(88, 5).Value = 2693
$ws.Cells.Item(88, 7).Value = 3
$ws.Cells.Item(88, 8).Value = 756
$ws.Cells.Item(89, 1).Value = 'Birmania'
$ws.Cells.Item(89, 2).Value = 17794
$ws.Cells.Item(89, 3).Value = 1291
$ws.Cells.Item(89, 4).Value = 5195
$ws.Cells.Item(89, 5).Value = 12187
$ws.Cells.Item(89, 7).Value = 41
$ws.Cells.Item(89, 8).Value = 412
$ws.Cells.Item(90, 1).Value = 'Croacia'
$ws.Cells.Item(90, 2).Value = 17659
$ws.Cells.Item(90, 3).Value = 258
$ws.Cells.Item(90, 4).Value = 15849
$ws.Cells.Item(90, 5).Value = 1512
$ws.Cells.Item(90, 7).Value = 5
$ws.Cells.Item(90, 8).Value = 298
$ws.Cells.Item(91, 1).Value = 'Madagascar'
$ws.Cells.Item(91, 2).Value = 16558
$ws.Cells.Item(91, 3).Value = 29
$ws.Cells.Item(91, 4).Value = 15486
$ws.Cells.Item(91, 5).Value = 840
$ws.Cells.Item(91, 8).Value = 232
$ws.Cells.Item(92, 2).Value = 15094
$ws.Cells.Item(92, 3).Value = 26
$ws.Cells.Item(92, 4).Value = 12805
$ws.Cells.Item(92, 5).Value = 1977
$ws.Cells.Item(95, 2).Value = 14421
$ws.Cells.Item(95, 3).Value = 59
$ws.Cells.Item(95, 5).Value = 2956
$ws.Cells.Item(103, 2).Value = 10754
$ws.Cells.Item(103, 3).Value = 19
$ws.Cells.Item(103, 4).Value = 10098
$ws.Cells.Item(103, 5).Value = 590
$ws.Cells.Item(110, 1).Value = 'Uganda'
$ws.Cells.Item(110, 2).Value = 8808
$ws.Cells.Item(110, 3).Value = 146
$ws.Cells.Item(110, 4).Value = 4736
$ws.Cells.Item(110, 5).Value = 3991
$ws.Cells.Item(110, 7).Value = 2
$ws.Cells.Item(110, 8).Value = 81
$ws.Cells.Item(111, 1).Value = 'Luxemburgo'
$ws.Cells.Item(111, 4).Value = 7428
$ws.Cells.Item(111, 5).Value = 1244
$ws.Cells.Item(111, 8).Value = 125
$ws.Cells.Item(112, 1).Value = 'Gabon'
$ws.Cells.Item(112, 2).Value = 8797
$ws.Cells.Item(112, 4).Value = 8067
$ws.Cells.Item(112, 5).Value = 676
$ws.Cells.Item(112, 8).Value = 54
$ws.Cells.Item(148, 2).Value = 2921
$ws.Cells.Item(148, 3).Value = 49
$ws.Cells.Item(148, 4).Value = 2277
$ws.Cells.Item(148, 5).Value = 634
$ws.Cells.Item(179, 2).Value = 474
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(179, 4).Value = 435
$ws.Cells.Item(179, 5).Value = 39
